$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: name/username edited to "fff" (password & e-mail stay "asd")
$ws.Range("A2").Value = "fff"
$ws.Range("B2").Value = "fff"

# New user row 4: Jose Manuel Martinez Peñaranda
$ws.Range("A4").Value = "José Manuel Martínez Peñaranda"
$ws.Range("B4").Value = "Rex117"
# Force the password / numeric-looking id to stay text, then drop the
# quote-prefix formatting it picks up so the cell keeps the default style
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "12345678"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").Value = "josemmp14@hotmail.com"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 3

# New user row 5
$ws.Range("A5").Value = "eee"
$ws.Range("B5").Value = "qqq"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "111"
$ws.Range("C5").ClearFormats()
$ws.Range("D5").Value = "sisoy"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 4
